# "Add files via upload" — add a new "Mehrfach" (multi-select) option row for
# security/access accessories on the Brix_Schiebe (sliding gate) catalog sheet,
# and move the active tab/selection from Zub_Zaun to Brix_Schiebe.

$wb = $excel.ActiveWorkbook

# --- Brix_Schiebe: insert the new "Sicherheits- & Zutritts-Zubehör" row ----
$wsSchiebe = $wb.Worksheets.Item("Brix_Schiebe")

# Push the existing Farbe/Antrieb/Preis rows (old rows 4-6) down by one.
$wsSchiebe.Rows.Item(4).Insert()

$wsSchiebe.Range("A4").Value = "Mehrfach"
$wsSchiebe.Range("B4").Value = "Sicherheits- & Zutritts-Zubehör"
$wsSchiebe.Range("C4").Value = "P_Zub"
$wsSchiebe.Range("D4").Value = "Lichtschranke:145, Blinkleuchte:95, Schlüsseltaster:120, GSM-Modul:350, Codetaster:180"
$wsSchiebe.Range("E4").Value = "Wird addiert"

# --- Move the active tab / selection to Brix_Schiebe ----------------------
# (previously Zub_Zaun, cell D2, was the selected/active sheet)
$wsSchiebe.Range("B13").Select()
